# Applies the "Add data for 2021-12-07" update to the carjacking-by-neighborhood-by-month
# workbook: renames the "through Nov 28" header/sheet name to "through Nov 29" and
# bumps/creates a handful of monthly count cells on the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and the header cell text (B1) that both describe the
# "as of" date for the current (November 2021) column.
$ws.Name = "Through 2021-11-29"
$ws.Range("B1").Value = "November 2021 (through November 29)"

# Cell-by-cell count updates / additions (row => column => new value).
$updates = @{
    3  = @{ "B" = 6;  "AT" = 7 }
    4  = @{ "M" = 13; "AI" = 8; "AT" = 12; "BE" = 5; "BP" = 4 }
    5  = @{ "M" = 7 }
    12 = @{ "B" = 5 }
    16 = @{ "M" = 6 }
    17 = @{ "AI" = 2 }
    19 = @{ "BE" = 1 }
    21 = @{ "AT" = 1 }
    24 = @{ "BP" = 1 }
    37 = @{ "M" = 1 }
    43 = @{ "B" = 1 }
    47 = @{ "AI" = 1 }
    48 = @{ "AI" = 2 }
    49 = @{ "M" = 3 }
    52 = @{ "AT" = 3 }
    61 = @{ "B" = 1 }
    66 = @{ "BP" = 1 }
    84 = @{ "BE" = 1 }
    97 = @{ "M" = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
